# Auto-generated script to apply scheduled market-data refresh to Gungnir_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Cells.Item(17, 8).Value = 56780.445
$ws.Cells.Item(17, 10).Value = 56780.445
$ws.Cells.Item(17, 12).Value = 170341.335
$ws.Cells.Item(17, 14).Value = -170677.335
# Row 40
$ws.Cells.Item(40, 8).Value = 2405812.2
$ws.Cells.Item(40, 9).Value = 5209813.5
$ws.Cells.Item(40, 10).Value = 2382.8572
$ws.Cells.Item(40, 11).Value = 5209813.5
$ws.Cells.Item(40, 12).Value = 2382.8572
$ws.Cells.Item(40, 13).Value = -5209638.5
$ws.Cells.Item(40, 14).Value = -2732.8572
# Row 62
$ws.Cells.Item(62, 8).Value = 12509371
$ws.Cells.Item(62, 9).Value = 31265276
$ws.Cells.Item(62, 10).Value = 5433.8335
$ws.Cells.Item(62, 11).Value = 31265276
$ws.Cells.Item(62, 12).Value = 5433.8335
$ws.Cells.Item(62, 13).Value = -31264652
$ws.Cells.Item(62, 14).Value = -6681.8335
# Row 65
$ws.Cells.Item(65, 8).Value = 12509371
$ws.Cells.Item(65, 9).Value = 31265276
$ws.Cells.Item(65, 10).Value = 5433.8335
$ws.Cells.Item(65, 11).Value = 156326380
$ws.Cells.Item(65, 12).Value = 27169.1675
$ws.Cells.Item(65, 13).Value = -156323260
$ws.Cells.Item(65, 14).Value = -33409.1675
# Row 98
$ws.Cells.Item(98, 8).Value = 83334250
$ws.Cells.Item(98, 9).Value = 89286450
$ws.Cells.Item(98, 11).Value = 89286450
$ws.Cells.Item(98, 13).Value = -89284952
# Row 122
$ws.Cells.Item(122, 8).Value = 83334250
$ws.Cells.Item(122, 9).Value = 89286450
$ws.Cells.Item(122, 11).Value = 267859350
$ws.Cells.Item(122, 13).Value = -267856900
# Row 125
$ws.Cells.Item(125, 8).Value = 1030
$ws.Cells.Item(125, 10).Value = 2000
$ws.Cells.Item(125, 12).Value = 18000
$ws.Cells.Item(125, 14).Value = -22920
# Row 138
$ws.Cells.Item(138, 8).Value = 3012.6118
$ws.Cells.Item(138, 9).Value = 1475.9688
$ws.Cells.Item(138, 10).Value = 3940.3962
$ws.Cells.Item(138, 11).Value = 4427.9064
$ws.Cells.Item(138, 12).Value = 11821.1886
$ws.Cells.Item(138, 13).Value = 712.0936000000002
$ws.Cells.Item(138, 14).Value = -22101.1886

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 6753.42
$ws.Cells.Item(32, 9).Value = 6653.0103
$ws.Cells.Item(32, 10).Value = 10000
$ws.Cells.Item(32, 11).Value = 6653.0103
$ws.Cells.Item(32, 12).Value = 10000
$ws.Cells.Item(32, 13).Value = -6366.0103
$ws.Cells.Item(32, 14).Value = -10574
# Row 45
$ws.Cells.Item(45, 8).Value = 15922191
$ws.Cells.Item(45, 9).Value = 18575472
$ws.Cells.Item(45, 10).Value = 2499.3333
$ws.Cells.Item(45, 11).Value = 18575472
$ws.Cells.Item(45, 12).Value = 2499.3333
$ws.Cells.Item(45, 13).Value = -18575095
$ws.Cells.Item(45, 14).Value = -3253.3333
# Row 61
$ws.Cells.Item(61, 8).Value = 1697.7819
$ws.Cells.Item(61, 9).Value = 1696.7
$ws.Cells.Item(61, 10).Value = 1699.08
$ws.Cells.Item(61, 11).Value = 1696.7
$ws.Cells.Item(61, 12).Value = 1699.08
$ws.Cells.Item(61, 13).Value = -1484.7
$ws.Cells.Item(61, 14).Value = -2123.08
# Row 88
$ws.Cells.Item(88, 8).Value = 2320.5186
$ws.Cells.Item(88, 9).Value = 2466.5
$ws.Cells.Item(88, 10).Value = 2234.647
$ws.Cells.Item(88, 11).Value = 2466.5
$ws.Cells.Item(88, 12).Value = 2234.647
$ws.Cells.Item(88, 13).Value = -2060.5
$ws.Cells.Item(88, 14).Value = -3046.647
# Row 91
$ws.Cells.Item(91, 8).Value = 2320.5186
$ws.Cells.Item(91, 9).Value = 2466.5
$ws.Cells.Item(91, 10).Value = 2234.647
$ws.Cells.Item(91, 11).Value = 2466.5
$ws.Cells.Item(91, 12).Value = 2234.647
$ws.Cells.Item(91, 13).Value = -1062.5
$ws.Cells.Item(91, 14).Value = -5042.647
# Row 122
$ws.Cells.Item(122, 8).Value = 2189.3809
$ws.Cells.Item(122, 9).Value = 2189.3809
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 6568.1427
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -4118.1427
$ws.Cells.Item(122, 14).ClearContents()
# Row 136
$ws.Cells.Item(136, 8).Value = 1697.7819
$ws.Cells.Item(136, 9).Value = 1696.7
$ws.Cells.Item(136, 10).Value = 1699.08
$ws.Cells.Item(136, 11).Value = 5090.1
$ws.Cells.Item(136, 12).Value = 5097.24
$ws.Cells.Item(136, 13).Value = -2540.1
$ws.Cells.Item(136, 14).Value = -10197.24

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Cells.Item(99, 8).Value = 2071.55
$ws.Cells.Item(99, 9).Value = 1250
$ws.Cells.Item(99, 10).Value = 2276.9375
$ws.Cells.Item(99, 11).Value = 1250
$ws.Cells.Item(99, 12).Value = 2276.9375
$ws.Cells.Item(99, 13).Value = 248
$ws.Cells.Item(99, 14).Value = -5272.9375
# Row 107
$ws.Cells.Item(107, 8).Value = 21740056
$ws.Cells.Item(107, 9).Value = 26316484
$ws.Cells.Item(107, 10).Value = 2015
$ws.Cells.Item(107, 11).Value = 26316484
$ws.Cells.Item(107, 12).Value = 2015
$ws.Cells.Item(107, 13).Value = -26314564
$ws.Cells.Item(107, 14).Value = -5855
# Row 122
$ws.Cells.Item(122, 8).Value = 40613.5
$ws.Cells.Item(122, 10).Value = 40613.5
$ws.Cells.Item(122, 12).Value = 40613.5
$ws.Cells.Item(122, 14).Value = -50413.5

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Cells.Item(7, 8).Value = 14740.292
$ws.Cells.Item(7, 9).Value = 38.272728
$ws.Cells.Item(7, 10).Value = 27180.46
$ws.Cells.Item(7, 11).Value = 38.272728
$ws.Cells.Item(7, 12).Value = 27180.46
$ws.Cells.Item(7, 13).Value = 74.727272
$ws.Cells.Item(7, 14).Value = -27406.46
# Row 22
$ws.Cells.Item(22, 8).Value = 10195.2
$ws.Cells.Item(22, 9).Value = 193.14285
$ws.Cells.Item(22, 11).Value = 193.14285
$ws.Cells.Item(22, 13).Value = 156.85715
# Row 86
$ws.Cells.Item(86, 8).Value = 22752440
$ws.Cells.Item(86, 9).Value = 35716390
$ws.Cells.Item(86, 10).Value = 65524.75
$ws.Cells.Item(86, 11).Value = 35716390
$ws.Cells.Item(86, 12).Value = 65524.75
$ws.Cells.Item(86, 13).Value = -35715267
$ws.Cells.Item(86, 14).Value = -67770.75
# Row 89
$ws.Cells.Item(89, 8).Value = 22752440
$ws.Cells.Item(89, 9).Value = 35716390
$ws.Cells.Item(89, 10).Value = 65524.75
$ws.Cells.Item(89, 11).Value = 178581950
$ws.Cells.Item(89, 12).Value = 327623.75
$ws.Cells.Item(89, 13).Value = -178576334
$ws.Cells.Item(89, 14).Value = -338855.75
# Row 107
$ws.Cells.Item(107, 8).Value = 525.13043
$ws.Cells.Item(107, 9).Value = 346.76923
$ws.Cells.Item(107, 10).Value = 757
$ws.Cells.Item(107, 11).Value = 346.76923
$ws.Cells.Item(107, 12).Value = 757
$ws.Cells.Item(107, 13).Value = 1573.23077
$ws.Cells.Item(107, 14).Value = -4597
# Row 134
$ws.Cells.Item(134, 8).Value = 1721.2222
$ws.Cells.Item(134, 9).Value = 1892.069
$ws.Cells.Item(134, 11).Value = 5676.207
$ws.Cells.Item(134, 13).Value = -3141.207
# Row 141
$ws.Cells.Item(141, 8).Value = 56304.8
$ws.Cells.Item(141, 10).Value = 63811.766
$ws.Cells.Item(141, 12).Value = 63811.766
$ws.Cells.Item(141, 14).Value = -74171.766

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Cells.Item(12, 8).Value = 67383.266
$ws.Cells.Item(12, 9).Value = 10.666667
$ws.Cells.Item(12, 11).Value = 32.000001
$ws.Cells.Item(12, 13).Value = 140.999999
# Row 98
$ws.Cells.Item(98, 8).Value = 624.1429000000001
$ws.Cells.Item(98, 9).Value = 0
$ws.Cells.Item(98, 10).Value = 624.1429000000001
$ws.Cells.Item(98, 11).Value = 0
$ws.Cells.Item(98, 12).Value = 1872.4287
$ws.Cells.Item(98, 14).Value = -4868.4287
$ws.Cells.Item(98, 13).ClearContents()
# Row 131
$ws.Cells.Item(131, 8).Value = 699
$ws.Cells.Item(131, 9).Value = 327.375
$ws.Cells.Item(131, 10).Value = 731.31525
$ws.Cells.Item(131, 11).Value = 982.125
$ws.Cells.Item(131, 12).Value = 2193.94575
$ws.Cells.Item(131, 13).Value = 4057.875
$ws.Cells.Item(131, 14).Value = -12273.94575

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Cells.Item(102, 8).Value = 1190.2222
$ws.Cells.Item(102, 9).Value = 1068.6666
$ws.Cells.Item(102, 10).Value = 1433.3334
$ws.Cells.Item(102, 11).Value = 1068.6666
$ws.Cells.Item(102, 12).Value = 1433.3334
$ws.Cells.Item(102, 13).Value = 553.3334
$ws.Cells.Item(102, 14).Value = -4677.3334
# Row 112
$ws.Cells.Item(112, 8).Value = 37041
$ws.Cells.Item(112, 10).Value = 37041
$ws.Cells.Item(112, 12).Value = 37041
$ws.Cells.Item(112, 14).Value = -39257
# Row 126
$ws.Cells.Item(126, 8).Value = 3526
$ws.Cells.Item(126, 9).Value = 6000
$ws.Cells.Item(126, 10).Value = 3062.125
$ws.Cells.Item(126, 11).Value = 18000
$ws.Cells.Item(126, 12).Value = 9186.375
$ws.Cells.Item(126, 13).Value = -15530
$ws.Cells.Item(126, 14).Value = -14126.375
# Row 132
$ws.Cells.Item(132, 8).Value = 8505.883
$ws.Cells.Item(132, 9).Value = 1511.3334
$ws.Cells.Item(132, 10).Value = 16374.75
$ws.Cells.Item(132, 11).Value = 4534.0002
$ws.Cells.Item(132, 12).Value = 49124.25
$ws.Cells.Item(132, 13).Value = -2004.0002
$ws.Cells.Item(132, 14).Value = -54184.25

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Cells.Item(16, 8).Value = 56392308
$ws.Cells.Item(16, 9).Value = 4762755.5
$ws.Cells.Item(16, 10).Value = 250003120
$ws.Cells.Item(16, 11).Value = 4762755.5
$ws.Cells.Item(16, 12).Value = 250003120
$ws.Cells.Item(16, 13).Value = -4762585.5
$ws.Cells.Item(16, 14).Value = -250003460
# Row 43
$ws.Cells.Item(43, 8).Value = 4000
$ws.Cells.Item(43, 10).Value = 4000
$ws.Cells.Item(43, 12).Value = 4000
$ws.Cells.Item(43, 14).Value = -4386
# Row 55
$ws.Cells.Item(55, 8).Value = 34486544
$ws.Cells.Item(55, 9).Value = 5562.5264
$ws.Cells.Item(55, 10).Value = 100000400
$ws.Cells.Item(55, 11).Value = 5562.5264
$ws.Cells.Item(55, 12).Value = 100000400
$ws.Cells.Item(55, 13).Value = -5389.5264
$ws.Cells.Item(55, 14).Value = -100000746
# Row 122
$ws.Cells.Item(122, 8).Value = 10452
$ws.Cells.Item(122, 9).Value = 12605.583
$ws.Cells.Item(122, 10).Value = 3991.25
$ws.Cells.Item(122, 11).Value = 37816.749
$ws.Cells.Item(122, 12).Value = 11973.75
$ws.Cells.Item(122, 13).Value = -35366.749
$ws.Cells.Item(122, 14).Value = -16873.75
# Row 136
$ws.Cells.Item(136, 8).Value = 2933.2727
$ws.Cells.Item(136, 9).Value = 2358.3513
$ws.Cells.Item(136, 10).Value = 5972.143
$ws.Cells.Item(136, 11).Value = 7075.053899999999
$ws.Cells.Item(136, 12).Value = 17916.429
$ws.Cells.Item(136, 13).Value = -4525.053899999999
$ws.Cells.Item(136, 14).Value = -23016.429

$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Cells.Item(100, 8).Value = 695.4375
$ws.Cells.Item(100, 9).Value = 792.6667
$ws.Cells.Item(100, 10).Value = 403.75
$ws.Cells.Item(100, 11).Value = 1585.3334
$ws.Cells.Item(100, 12).Value = 807.5
$ws.Cells.Item(100, 13).Value = -1044.3334
$ws.Cells.Item(100, 14).Value = -1889.5
# Row 123
$ws.Cells.Item(123, 8).Value = 17347
$ws.Cells.Item(123, 10).Value = 17347
$ws.Cells.Item(123, 12).Value = 17347
$ws.Cells.Item(123, 14).Value = -27147
# Row 126
$ws.Cells.Item(126, 8).Value = 1237.3125
$ws.Cells.Item(126, 9).Value = 928.9167
$ws.Cells.Item(126, 10).Value = 2162.5
$ws.Cells.Item(126, 11).Value = 2786.7501
$ws.Cells.Item(126, 12).Value = 6487.5
$ws.Cells.Item(126, 13).Value = -316.7501000000002
$ws.Cells.Item(126, 14).Value = -11427.5
